$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 2598
$ws.Range("D6").Value = 2621
$ws.Range("E6").Value = 2650
$ws.Range("F6").Value = 2650
$ws.Range("G6").Value = 2698
